# Bitacora_tareas.xlsx - add new "Diagnostico" log entries (rows 25-27 on
# "Log") documenting the classification-error exceptions, plus the matching
# version-history row (row 5 on "Versiones").

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Log": append rows 25, 26 and 27
# ---------------------------------------------------------------------
$log = $wb.Worksheets.Item("Log")

$log.Cells.Item(25, 1).Value = "27/02/2025"
$log.Cells.Item(25, 2).Value = "17:40"
$log.Cells.Item(25, 3).Value = "Excepción errores: Comisiones Bancarias / Gastos Bancarios"
$log.Cells.Item(25, 4).Value = "Si la categoría es Comisiones Bancarias y la cuenta contable es Gastos Bancarios, se considera consistente y no entra en el log de errores de clasificación (aunque la descripción no contenga esas palabras)."
$log.Cells.Item(25, 5).Value = "Diagnostico"

$log.Cells.Item(26, 1).Value = "27/02/2025"
$log.Cells.Item(26, 2).Value = "17:50"
$log.Cells.Item(26, 3).Value = "Excepción errores: Impuestos / MercadoPago y Impuestos / Transferencia Morba"
$log.Cells.Item(26, 4).Value = "Si la categoría es Impuestos y la cuenta contable es MercadoPago o Transferencia Morba, se considera consistente y no entra en el log de errores de clasificación, aunque la descripción no contenga esas palabras."
$log.Cells.Item(26, 5).Value = "Diagnostico"

$log.Cells.Item(27, 1).Value = "27/02/2025"
$log.Cells.Item(27, 2).Value = "18:00"
$log.Cells.Item(27, 3).Value = "Excepción errores: Alquileres y Servicios / Alquiler"
$log.Cells.Item(27, 4).Value = "Si la categoría es Alquiler (mostrada como Alquileres y Servicios) y la cuenta contable es Alquiler, se considera consistente y no entra en el log de errores de clasificación."
$log.Cells.Item(27, 5).Value = "Diagnostico"

# ---------------------------------------------------------------------
# Sheet "Versiones": append row 5 (version 1.3)
# ---------------------------------------------------------------------
$versiones = $wb.Worksheets.Item("Versiones")

# "1.3" looks numeric, so force it to stay text (matching the existing
# Version column entries "1.0"/"1.1"/"1.2") the same way a user would in
# the Excel UI: a leading apostrophe.
$versiones.Cells.Item(5, 1).Value = "'1.3"
$versiones.Cells.Item(5, 2).Value = "27/02/2025"
$versiones.Cells.Item(5, 3).Value = "Errores de clasificación (solapa Errores), edición desde modal, editado/editado_detalle; excepciones: Comisiones Bancarias/Gastos Bancarios, Impuestos/MercadoPago y Transferencia Morba, Alquiler/Alquiler"

# ---------------------------------------------------------------------
# Mirror Excel's own "ignore the number-stored-as-text warning" action for
# the full (now larger) tables, same as the green-triangle indicator the
# user dismisses over the whole range. No-op if unsupported by the host.
# ---------------------------------------------------------------------
try { $log.Range("A1:E27").Errors.Item(9).Ignore = $true } catch { }
try { $versiones.Range("A1:C5").Errors.Item(9).Ignore = $true } catch { }

Write-Host "Log dimension:" $log.UsedRange.Address()
Write-Host "Versiones dimension:" $versiones.UsedRange.Address()
